# Applies the scheduled-runner value updates to the Spriggan_Profits sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1109.9131
$ws.Range("I41").Value = 1368.2222
$ws.Range("K41").Value = 1368.2222
$ws.Range("M41").Value = -928.2221999999999
$ws.Range("H112").Value = 98663.91
$ws.Range("I112").Value = 252336
$ws.Range("J112").Value = 62505.766
$ws.Range("K112").Value = 757008
$ws.Range("L112").Value = 187517.298
$ws.Range("M112").Value = -755900
$ws.Range("N112").Value = -189733.298
$ws.Range("H124").Value = 59999
$ws.Range("J124").Value = 59999
$ws.Range("L124").Value = 59999
$ws.Range("N124").Value = -69819
$ws.Range("H132").Value = 2536
$ws.Range("I132").Value = 2536
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7608
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5078
$ws.Range("N132").Value = $null
$ws.Range("H135").Value = 21739678
$ws.Range("I135").Value = 21739678
$ws.Range("K135").Value = 195657102
$ws.Range("M135").Value = -195654567
$ws.Range("H137").Value = 2407.457
$ws.Range("I137").Value = 2119.9583
$ws.Range("K137").Value = 6359.874899999999
$ws.Range("M137").Value = -3809.874899999999
$ws.Range("H141").Value = 712.65625
$ws.Range("I141").Value = 755.5172
$ws.Range("J141").Value = 298.33334
$ws.Range("K141").Value = 2266.5516
$ws.Range("L141").Value = 895.0000200000001
$ws.Range("M141").Value = 2913.4484
$ws.Range("N141").Value = -11255.00002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2468.1277
$ws.Range("I32").Value = 2361.6943
$ws.Range("J32").Value = 2816.4546
$ws.Range("K32").Value = 2361.6943
$ws.Range("L32").Value = 2816.4546
$ws.Range("M32").Value = -2074.6943
$ws.Range("N32").Value = -3390.4546
$ws.Range("H88").Value = 144397.58
$ws.Range("J88").Value = 1836.6
$ws.Range("L88").Value = 1836.6
$ws.Range("N88").Value = -2648.6
$ws.Range("H91").Value = 144397.58
$ws.Range("J91").Value = 1836.6
$ws.Range("L91").Value = 1836.6
$ws.Range("N91").Value = -4644.6
$ws.Range("H110").Value = 54962.26
$ws.Range("I110").Value = 68472.266
$ws.Range("J110").Value = 4299.75
$ws.Range("K110").Value = 68472.266
$ws.Range("L110").Value = 4299.75
$ws.Range("M110").Value = -66427.266
$ws.Range("N110").Value = -8389.75
$ws.Range("H125").Value = 50500
$ws.Range("J125").Value = 50500
$ws.Range("L125").Value = 50500
$ws.Range("N125").Value = -60340

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1882
$ws.Range("I105").Value = 1882
$ws.Range("K105").Value = 1882
$ws.Range("M105").Value = -135
$ws.Range("H122").Value = 39700
$ws.Range("J122").Value = 39700
$ws.Range("L122").Value = 39700
$ws.Range("N122").Value = -49500
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8330.257
$ws.Range("I31").Value = 6191.4116
$ws.Range("J31").Value = 10350.277
$ws.Range("K31").Value = 6191.4116
$ws.Range("L31").Value = 10350.277
$ws.Range("M31").Value = -5896.4116
$ws.Range("N31").Value = -10940.277
$ws.Range("H34").Value = 8330.257
$ws.Range("I34").Value = 6191.4116
$ws.Range("J34").Value = 10350.277
$ws.Range("K34").Value = 6191.4116
$ws.Range("L34").Value = 10350.277
$ws.Range("M34").Value = -5989.4116
$ws.Range("N34").Value = -10754.277
$ws.Range("H58").Value = 16671083
$ws.Range("I58").Value = 22732728
$ws.Range("J58").Value = 1558
$ws.Range("K58").Value = 22732728
$ws.Range("L58").Value = 1558
$ws.Range("M58").Value = -22732525
$ws.Range("N58").Value = -1964
$ws.Range("H122").Value = 2937.25
$ws.Range("I122").Value = 2937.25
$ws.Range("K122").Value = 8811.75
$ws.Range("M122").Value = -6361.75
$ws.Range("H132").Value = 20409828
$ws.Range("I132").Value = 21278272
$ws.Range("K132").Value = 63834816
$ws.Range("M132").Value = -63832286
$ws.Range("H136").Value = 16671083
$ws.Range("I136").Value = 22732728
$ws.Range("J136").Value = 1558
$ws.Range("K136").Value = 68198184
$ws.Range("L136").Value = 4674
$ws.Range("M136").Value = -68195634
$ws.Range("N136").Value = -9774

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 274.75
$ws.Range("I8").Value = 274.75
$ws.Range("K8").Value = 824.25
$ws.Range("M8").Value = -685.25
$ws.Range("H137").Value = 16668158
$ws.Range("I137").Value = 20001388
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 60004164
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -59999064
$ws.Range("N137").Value = -16200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 64999.668
$ws.Range("J124").Value = 64999.668
$ws.Range("L124").Value = 64999.668
$ws.Range("N124").Value = -74819.66800000001
$ws.Range("H126").Value = 7546.3335
$ws.Range("I126").Value = 7805.636
$ws.Range("J126").Value = 4694
$ws.Range("K126").Value = 23416.908
$ws.Range("L126").Value = 14082
$ws.Range("M126").Value = -20946.908
$ws.Range("N126").Value = -19022
$ws.Range("H132").Value = 20835026
$ws.Range("I132").Value = 31251938
$ws.Range("K132").Value = 93755814
$ws.Range("M132").Value = -93753284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1235.7142
$ws.Range("I107").Value = 1100.8
$ws.Range("J107").Value = 1573
$ws.Range("K107").Value = 3302.4
$ws.Range("L107").Value = 4719
$ws.Range("M107").Value = -1382.4
$ws.Range("N107").Value = -8559
$ws.Range("H126").Value = 1165.7142
$ws.Range("I126").Value = 998.8889
$ws.Range("K126").Value = 2996.6667
$ws.Range("M126").Value = -526.6667000000002
$ws.Range("H132").Value = 10207440
$ws.Range("I132").Value = 11364942
$ws.Range("K132").Value = 34094826
$ws.Range("M132").Value = -34092296
$ws.Range("H136").Value = 17243792
$ws.Range("I136").Value = 17859606
$ws.Range("K136").Value = 53578818
$ws.Range("M136").Value = -53576268
